$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.508.64'
$ws.Range('D3').Value = '3.587.53'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.12%  '
$ws.Range('D7').Value = '3.588.29'
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.495'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.125'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.20'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.392'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').Value = '4.187.08'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000186'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.22%  '
$ws.Range('D16').Value = '3.579.30'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.117'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '64.678.89'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.05'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.60'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '394.69'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.587'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.16%  '
$ws.Range('D24').Value = '3.729.20'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.33'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +26.73%  '
$ws.Range('E30').Value = '  +3.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.68'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.34%  '
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('D33').Value = '3.583.41'
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.50'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  -0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.60'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.08'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '171.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0831'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.846'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.33'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.52'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.79%  '
$ws.Range('D50').Value = '2.458.60'
$ws.Range('E50').Value = '  -1.26%  '
$ws.Range('E51').Value = '  +2.10%  '
